# Adds the first benchmark timings for bf2 on the GTX1070 (rows DE/VT/ME,
# i.e. rows 16-18 under "bf2-aos-sh" and rows 26-28 under "bf2-soa-nosh"),
# plus marks the "rome" / timeout runs (NV rows 9/19/29) as ">1h" (SPEEDUP
# column) / "N/A" (THROUGHPUT columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- bf2-aos-nosh block (rows 5-13), NV row (row 9): timed out ---
$ws.Range("I9").Value = ">1h"
$ws.Range("AD9:AH9").Value = "N/A"

# --- bf2-aos-sh block (rows 15-23) ---
# DE (row16), VT (row17), ME (row18): first GTX1070 timings
$ws.Range("I16").Value = 101.76300000000001
$ws.Range("AD16").Value = 57.79

$ws.Range("I17").Value = 397.11799999999999
$ws.Range("AD17").Value = 52.55

$ws.Range("I18").Value = 1564.66
$ws.Range("AD18").Value = 52.92

# NV (row19): timed out
$ws.Range("I19").Value = ">1h"
$ws.Range("AD19:AH19").Value = "N/A"

# --- bf2-soa-nosh block (rows 25-33) ---
# DE (row26), VT (row27), ME (row28): first GTX1070 timings
$ws.Range("I26").Value = 64.331999999999994
$ws.Range("AD26").Value = 91.41

$ws.Range("I27").Value = 368.23700000000002
$ws.Range("AD27").Value = 56.67

$ws.Range("I28").Value = 1448.54
$ws.Range("AD28").Value = 57.16

# NV (row29): timed out
$ws.Range("I29").Value = ">1h"
$ws.Range("AD29:AH29").Value = "N/A"

# --- view state: last selection before save ---
$ws.Range("AE31").Select()
